$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '67.276.29'
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").Value = '3.943.49'
$ws.Range("E3").Value = '  +4.16%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '470.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.16%  '

$ws.Range("E7").Value = '  +1.21%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.735'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("E10").Value = '  +9.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000343'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.47%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '4.555.68'
$ws.Range("E13").Value = '  +3.59%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("E15").Value = '  +2.23%  '

$ws.Range("D16").Value = '3.963.33'
$ws.Range("E16").Value = '  +5.38%  '

$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("E19").Value = '  +1.94%  '

$ws.Range("D20").Value = '67.512.96'
$ws.Range("E20").Value = '  +0.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '437.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.91%  '

$ws.Range("E22").Value = '  +4.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '719.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.44%  '

$ws.Range("E32").Value = '  +2.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.152'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.22%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.88'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.07%  '

$ws.Range("E36").Value = '  +20.21%  '

$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("E38").Value = '  -4.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0478'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.18%  '

$ws.Range("E42").Value = '  +0.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.338'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.82%  '

$ws.Range("E44").Value = '  -0.30%  '

$ws.Range("E45").Value = '  -7.45%  '

$ws.Range("E46").Value = '  +4.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '147.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.86%  '

